$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 126:127. This shifts the former rows 126 and 127
# (and everything below) down to 128 and 129.
$ws.Rows("126:127").Insert()

# Copy the formatting (style) of the index column from a known-good row (125)
# onto the two newly inserted index cells so they keep style s="1" (bold, boxed,
# centered) just like every other row-index cell in column A.
$ws.Cells.Item(125, 1).Copy()
$ws.Cells.Item(126, 1).PasteSpecial(-4122)
$ws.Cells.Item(127, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 126 is the new "125th" match entry.
$ws.Cells.Item(126, 1).Value = 125
$ws.Cells.Item(126, 7).Value = -14.91601092197619
$ws.Cells.Item(126, 8).Value = 14.91601092197619

# Row 127 is the new "126th" match entry.
$ws.Cells.Item(127, 1).Value = 126
$ws.Cells.Item(127, 7).Value = -12.41739680548837
$ws.Cells.Item(127, 8).Value = 12.41739680548837

# The rows that were pushed down (previously 126 and 127) keep their original
# values but now represent entries 127 and 128 respectively.
$ws.Cells.Item(128, 1).Value = 127
$ws.Cells.Item(129, 1).Value = 128
$ws.Cells.Item(129, 5).Value = -16.31453075323677
$ws.Cells.Item(129, 8).Value = 16.31453075323677
